$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF") -- match formatting of the
# existing header cells (bold, bordered, centered) by copying H1's format.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Per-row values for columns I ("I0") and J ("IF") -- identical values per row
$values = @(8,8,8,8,8,8,8,8,8,8,10,8,8,6,7,9,9,9,8,8,8,8,8,6,8,8,11,8,8,10,8,8,8,8,8,8,8,6,8,8,8,8,8,8,8,8,8,8,8,8,8,8,7,8,8,8,8,7,8,7,6,5)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $v = $values[$i]
    $ws.Cells.Item($row, 9).Value = $v
    $ws.Cells.Item($row, 10).Value = $v
}
